# Trade #42 closed at 2026-02-17 12:47:47 - unknown UNKNOWN +0.000%
#
# Updates the rolled-up summary figures (Summary + Strategy Status sheets)
# and appends the newly-closed trade as row 43 on both the "All Trades" and
# "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.1                 # Current Capital
$summary.Range("B4").Value = 0.09                   # Total P&L $
$summary.Range("B5").Value = 0.04                   # Total P&L %
$summary.Range("B6").Value = 42                     # Total Trades
$summary.Range("B7").Value = 17                     # Winning Trades
$summary.Range("B9").Value = 40.48                  # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.1                   # Capital
$status.Range("D4").Value = 42                      # Trades
$status.Range("E4").Value = 0.09                    # P&L $
$status.Range("F4").Value = 0.1                     # P&L %
$status.Range("G4").Value = 40.48                   # Win Rate %

# ---------------------------------------------------------------------------
# Append the closed trade (row 43) to "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------------
function Add-TradeRow43($ws) {
    $ws.Cells.Item(43, 1).Value = 42                  # Trade #
    # Date/Time look like dates to the smart-parser, so force literal text
    # with a leading apostrophe (same as typing '2026-02-17 into Excel).
    $ws.Cells.Item(43, 2).Value = "'2026-02-17"        # Date
    $ws.Cells.Item(43, 3).Value = "'12:47:41"          # Time
    $ws.Cells.Item(43, 4).Value = "MarketMaking"       # Strategy
    $ws.Cells.Item(43, 5).Value = "UP"                 # Side
    $ws.Cells.Item(43, 6).Value = 0.87                 # Entry Price
    $ws.Cells.Item(43, 7).Value = 0.89                 # Exit Price
    $ws.Cells.Item(43, 8).Value = "CLOSED"             # Status
    $ws.Cells.Item(43, 9).Value = 2.2989               # P&L %
    $ws.Cells.Item(43, 10).Value = 0.02                # P&L $
    $ws.Cells.Item(43, 11).Value = 100.1               # Capital After
    $ws.Cells.Item(43, 12).Value = 0                   # Entry Slippage (bps)
    $ws.Cells.Item(43, 13).Value = 0                   # Exit Slippage (bps)
    $ws.Cells.Item(43, 14).Value = 0.6                 # Confidence
    $ws.Cells.Item(43, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(43, 16).Value = "early_exit"        # Exit Reason
    $ws.Cells.Item(43, 17).Value = 0.13                # Duration (min)
}

Add-TradeRow43 $wb.Worksheets.Item("All Trades")
Add-TradeRow43 $wb.Worksheets.Item("MarketMaking")
